$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F5").Value = -3
$ws.Range("F8").Value = -3
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = -3
